$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final ordered data (language, value), sorted descending by value,
# with "Swedish" and "Uzbek" removed.
$data = @(
    @("English", 22.16380015235244),
    @("Chinese", 15.7883606718513),
    @("Spanish", 7.089477018739032),
    @("Arabic", 5.115424486784671),
    @("Japanese", 4.836192817025105),
    @("German", 4.432330044704834),
    @("Portuguese", 3.597331326148079),
    @("Russian", 3.568624866346352),
    @("Malay-Indonesian", 2.987970667658099),
    @("French", 2.748710123664888),
    @("Italian", 2.27024335846713),
    @("Korean", 1.697874429719104),
    @("Persian", 1.566395521779735),
    @("Turkish", 1.519065747487474),
    @("Dutch", 1.293414106723533),
    @("Thai", 0.9535634281603422),
    @("Polish", 0.9009768793958561),
    @("Urdu", 0.8608897953690846),
    @("Vietnamese", 0.5347545762744657),
    @("Bengali", 0.526834940488706)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Remove the two now-extra trailing rows (previously rows 22 and 23).
$ws.Range("A22:B23").Delete()
